# Applies the "Updated cryptos list" price/volume refresh.
# D-column cells get NumberFormat="@" first so Excel stores the
# literal text (preserves trailing zeros / leading zeros / multi-dot
# "thousands" separators) instead of silently coercing to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.669.51"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.243.86"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.37"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "298.02"
$ws.Range("E6").Value = "  +11.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.618"
$ws.Range("E9").Value = "  +0.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.03"
$ws.Range("E10").Value = "  -4.11%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.98"
$ws.Range("E11").Value = "  +9.93%  "
$ws.Range("B12").Value = "Dogecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0930"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.08"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("E14").Value = "  -2.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.31"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.892"
$ws.Range("E16").Value = "  +1.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.580.89"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.260.54"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.730.02"
$ws.Range("E19").Value = "  -2.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.61"
$ws.Range("E20").Value = "  +10.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0000108"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.64"
$ws.Range("E22").Value = "  +25.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.52"
$ws.Range("E23").Value = "  +1.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").Value = "  -4.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "231.97"
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.46"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  +4.02%  "
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.91"
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.40"
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.27"
$ws.Range("E33").Value = "  -2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0905"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.69"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.54"
$ws.Range("E36").Value = "  +14.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.128"
$ws.Range("E37").Value = "  -1.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.79"
$ws.Range("E38").Value = "  +1.91%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0371"
$ws.Range("E39").Value = "  -3.89%  "
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.59"
$ws.Range("E41").Value = "  +0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.240"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.80"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.36"
$ws.Range("E44").Value = "  -7.05%  "
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.34"
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.56"
$ws.Range("E47").Value = "  -7.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.38"
$ws.Range("E48").Value = "  +7.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.40"
$ws.Range("E49").Value = "  +6.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.62"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0988"
$ws.Range("E51").Value = "  -2.02%  "
